$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) holds text values that look numeric (e.g. "565.63" or
# thousands-dotted "59.481.46"). Pre-formatting the cell as Text ("@") before
# assigning keeps Excel from silently converting these into numbers, matching
# the source data which is stored as plain text.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.481.46"
$ws.Range("E2").Value = "  +2.88%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.983.82"
$ws.Range("E3").Value = "  +1.28%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "565.63"
$ws.Range("E5").Value = "  +2.39%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.58"
$ws.Range("E6").Value = "  +4.01%  "

$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.521"
$ws.Range("E8").Value = "  +1.26%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.975.22"
$ws.Range("E9").Value = "  +1.16%  "

$ws.Range("E10").Value = "  +3.60%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.39"
$ws.Range("E11").Value = "  +12.06%  "

$ws.Range("E12").Value = "  +0.27%  "

$ws.Range("E13").Value = "  +4.25%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.75"
$ws.Range("E14").Value = "  +2.83%  "

$ws.Range("E15").Value = "  +0.08%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.474.12"
$ws.Range("E16").Value = "  +1.15%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.05"
$ws.Range("E17").Value = "  +1.40%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.978.66"
$ws.Range("E18").Value = "  +1.15%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "59.477.05"
$ws.Range("E19").Value = "  +2.89%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "436.54"
$ws.Range("E20").Value = "  +4.94%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.62"
$ws.Range("E21").Value = "  +2.00%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.721"
$ws.Range("E22").Value = "  +3.40%  "

$ws.Range("E23").Value = "  +0.21%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.31"
$ws.Range("E24").Value = "  -1.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "79.97"
$ws.Range("E25").Value = "  +1.25%  "

$ws.Range("E26").Value = "  -0.15%  "

$ws.Range("E27").Value = "  +10.02%  "

$ws.Range("E28").Value = "  -0.16%  "

$ws.Range("E29").Value = "  +2.35%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.75"
$ws.Range("E30").Value = "  +3.41%  "

$ws.Range("E31").Value = "  +1.14%  "

$ws.Range("E32").Value = "  +4.34%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.104"
$ws.Range("E33").Value = "  +8.52%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0775"
$ws.Range("E34").Value = "  +11.43%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.91"
$ws.Range("E35").Value = "  +4.14%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.979"
$ws.Range("E36").Value = "  +3.42%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.07"
$ws.Range("E37").Value = "  +0.59%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.62"
$ws.Range("E38").Value = "  +0.58%  "

$ws.Range("E39").Value = "  -3.53%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.79"
$ws.Range("E40").Value = "  +2.97%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "398.80"
$ws.Range("E41").Value = "  +4.53%  "

$ws.Range("E42").Value = "  +1.26%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.734.19"
$ws.Range("E43").Value = "  +1.01%  "

$ws.Range("E44").Value = "  -2.28%  "

$ws.Range("E45").Value = "  +5.84%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "35.23"
$ws.Range("E46").Value = "  +21.09%  "

$ws.Range("E47").Value = "  +0.00%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "122.13"
$ws.Range("E48").Value = "  -1.85%  "

$ws.Range("E49").Value = "  +1.92%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.01"
$ws.Range("E50").Value = "  +1.79%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.27"
$ws.Range("E51").Value = "  +1.73%  "
